$d = $word.ActiveDocument

# Table with the person temporarily covering the position:
# Apellido Paterno / Apellido Materno / Nombre(s) -> all become "EJEMPLO"
$d.Content.Find.Execute("SALVADOR", $true, $false, $false, $false, $false,
                         $true, 1, $false, "EJEMPLO", 2)
$d.Content.Find.Execute("JIMÉNEZ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "EJEMPLO", 2)
$d.Content.Find.Execute("ISIDRO NOÉ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "EJEMPLO", 2)

# "A PARTIR DEL:" date changes from 03 NOVIEMBRE to 02 FEBRERO
$d.Content.Find.Execute("03", $true, $false, $false, $false, $false,
                         $true, 1, $false, "02", 2)
$d.Content.Find.Execute("NOVIEMBRE", $true, $false, $false, $false, $false,
                         $true, 1, $false, "FEBRERO", 2)
